$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.165.33"
$ws.Range("E2").Value = "  +3.50%  "
$ws.Range("D3").Value = "2.439.10"
$ws.Range("E3").Value = "  +5.20%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'556.03"
$ws.Range("E5").Value = "  +2.42%  "
$ws.Range("D6").Value = "'139.04"
$ws.Range("E6").Value = "  +6.35%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'0.583"
$ws.Range("E8").Value = "  +1.17%  "
$ws.Range("D9").Value = "2.435.49"
$ws.Range("E9").Value = "  +5.14%  "
$ws.Range("D10").Value = "'0.105"
$ws.Range("E10").Value = "  +3.62%  "
$ws.Range("D11").Value = "'5.74"
$ws.Range("E11").Value = "  +3.87%  "
$ws.Range("D12").Value = "'0.151"
$ws.Range("E12").Value = "  +0.29%  "
$ws.Range("D13").Value = "'0.349"
$ws.Range("E13").Value = "  +5.10%  "
$ws.Range("D14").Value = "'26.12"
$ws.Range("E14").Value = "  +11.48%  "
$ws.Range("D15").Value = "2.872.14"
$ws.Range("E15").Value = "  +5.26%  "
$ws.Range("D16").Value = "62.047.07"
$ws.Range("E16").Value = "  +3.37%  "
$ws.Range("D17").Value = "'0.0000142"
$ws.Range("E17").Value = "  +7.77%  "
$ws.Range("D18").Value = "2.431.49"
$ws.Range("E18").Value = "  +5.19%  "
$ws.Range("D19").Value = "'11.20"
$ws.Range("E19").Value = "  +6.61%  "
$ws.Range("D20").Value = "'344.69"
$ws.Range("E20").Value = "  +10.11%  "
$ws.Range("D21").Value = "'4.19"
$ws.Range("E21").Value = "  +2.44%  "
$ws.Range("E22").Value = "  +3.05%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").Value = "'65.10"
$ws.Range("E24").Value = "  +2.14%  "
$ws.Range("E25").Value = "  +0.82%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("D27").Value = "'1.53"
$ws.Range("E27").Value = "  +13.70%  "
$ws.Range("D28").Value = "'8.23"
$ws.Range("E28").Value = "  +5.96%  "
$ws.Range("E29").Value = "  +12.57%  "
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "0.0₃0787"
$ws.Range("E30").Value = "  +8.25%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'1.80"
$ws.Range("E31").Value = "  +4.94%  "
$ws.Range("B32").Value = "Aptos"
$ws.Range("C32").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D32").Value = "'6.35"
$ws.Range("E32").Value = "  +8.25%  "
$ws.Range("B33").Value = "Monero"
$ws.Range("C33").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D33").Value = "'171.74"
$ws.Range("E33").Value = "  +0.52%  "
$ws.Range("D34").Value = "'1.44"
$ws.Range("E34").Value = "  +4.99%  "
$ws.Range("D35").Value = "'0.396"
$ws.Range("E35").Value = "  +4.53%  "
$ws.Range("D36").Value = "'18.59"
$ws.Range("E36").Value = "  +5.14%  "
$ws.Range("D37").Value = "'4.48"
$ws.Range("E37").Value = "  +11.60%  "
$ws.Range("D38").Value = "'368.46"
$ws.Range("E38").Value = "  +15.81%  "
$ws.Range("E40").Value = "  -0.11%  "
$ws.Range("D41").Value = "'1.70"
$ws.Range("E41").Value = "  +11.39%  "
$ws.Range("D42").Value = "'39.25"
$ws.Range("E42").Value = "  +3.73%  "
$ws.Range("D43").Value = "'146.55"
$ws.Range("E43").Value = "  +7.36%  "
$ws.Range("D44").Value = "'3.67"
$ws.Range("E44").Value = "  +6.45%  "
$ws.Range("D45").Value = "'20.65"
$ws.Range("E45").Value = "  +9.71%  "
$ws.Range("E46").Value = "  +1.50%  "
$ws.Range("D47").Value = "'0.589"
$ws.Range("E47").Value = "  +4.38%  "
$ws.Range("D48").Value = "'0.0518"
$ws.Range("E48").Value = "  +5.45%  "
$ws.Range("D49").Value = "'0.0222"
$ws.Range("E49").Value = "  +4.79%  "
$ws.Range("D50").Value = "'17.84"
$ws.Range("E50").Value = "  +6.28%  "
$ws.Range("D51").Value = "0.0₆0218"
$ws.Range("E51").Value = "  -2.70%  "
